$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.181459307670593
$ws.Range("B1").Value = 2.400239706039429
$ws.Range("C1").Value = 3.711603164672852
$ws.Range("D1").Value = 2.079583406448364
$ws.Range("E1").Value = 1.202973961830139
